# Weekly Update and cycle sort fix
#
# Adds four new leaderboard run rows (87-90) below the existing table and
# tidies up the formatting of what used to be the last row (86), which
# picks up the "normal" look shared by the rest of the data rows now that
# it is no longer the final row of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

$newRows = @(
    @{ Row = 87; A = "Yangi";  B = 10; C = "https://www.youtube.com/watch?v=d5_vxtrOUg8"; D = 0.23680555555555557; E = "Dehya"; F = "Kazuha";  G = "Bennett";   H = "Rosaria"; I = "Yae";    J = "Nahida";  K = "Tighnari"; L = "Lisa" },
    @{ Row = 88; A = "Yangi";  B = 10; C = "https://www.youtube.com/watch?v=Z4hsqVwV308"; D = 0.19513888888888889; E = "HuTao"; F = "Xingqiu"; G = "Mona";      H = "Dehya";   I = "Yae";    J = "Nahida";  K = "Tighnari"; L = "Lisa" },
    @{ Row = 89; A = "KennyL"; B = 10; C = "https://youtu.be/VIh3RUIfBOI";                D = 0.27569444444444446; E = "Razor"; F = "Xingqiu"; G = "Bennett";   H = "Nahida";  I = "Shenhe"; J = "Rosaria"; K = "Chongyun"; L = "Kazuha" },
    @{ Row = 90; A = "KennyL"; B = 10; C = "https://youtu.be/MDbqjALh8Po";                D = 0.20902777777777778; E = "Razor"; F = "Bennett"; G = "Xiangling"; H = "Jean";    I = "Shenhe"; J = "Rosaria"; K = "Chongyun"; L = "Kazuha" }
)

# ---------------------------------------------------------------------
# 1) Stamp the new rows with the same formatting already used by row 86
#    (borders / fonts / centred+wrapped text) before writing any values,
#    so the appended rows look like a natural continuation of the table.
# ---------------------------------------------------------------------
$ws.Range("A86:L86").Copy()
$ws.Range("A87:L87").PasteSpecial($xlPasteFormats)
$ws.Range("A88:L88").PasteSpecial($xlPasteFormats)
$ws.Range("A89:L89").PasteSpecial($xlPasteFormats)
$ws.Range("A90:L90").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Column C (Video) uses its own centred style (matching the rest of the
# Video column) rather than row 86's style.
$ws.Range("C85").Copy()
$ws.Range("C87").PasteSpecial($xlPasteFormats)
$ws.Range("C88").PasteSpecial($xlPasteFormats)
$ws.Range("C89").PasteSpecial($xlPasteFormats)
$ws.Range("C90").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 4).NumberFormat = "h:mm"
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
}

# ---------------------------------------------------------------------
# 2) Row 86 is no longer the last row of the table, so bring its look
#    back in line with the rest of the data rows (e.g. row 85), leaving
#    its Video cell (column C) untouched.
# ---------------------------------------------------------------------
$ws.Range("A85:B85").Copy()
$ws.Range("A86:B86").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("D85:L85").Copy()
$ws.Range("D86:L86").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Refresh the view so the newly added rows are in focus.
# ---------------------------------------------------------------------
$ws.Range("M90").Select()
